$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 2-44: update Price (D) and Volume(1h) (E) text values
$c = $ws.Range('D2'); $c.NumberFormat = "@"; $c.Value = '29.109.84'; $c.Style = "Normal"
$c = $ws.Range('E2'); $c.NumberFormat = "@"; $c.Value = '  -1.73%  '; $c.Style = "Normal"
$c = $ws.Range('D3'); $c.NumberFormat = "@"; $c.Value = '1.837.38'; $c.Style = "Normal"
$c = $ws.Range('E3'); $c.NumberFormat = "@"; $c.Value = '  -1.24%  '; $c.Style = "Normal"
$c = $ws.Range('D4'); $c.NumberFormat = "@"; $c.Value = '0.9991'; $c.Style = "Normal"
$c = $ws.Range('E4'); $c.NumberFormat = "@"; $c.Value = '  -0.03%  '; $c.Style = "Normal"
$c = $ws.Range('D5'); $c.NumberFormat = "@"; $c.Value = '240.61'; $c.Style = "Normal"
$c = $ws.Range('E5'); $c.NumberFormat = "@"; $c.Value = '  -1.74%  '; $c.Style = "Normal"
$c = $ws.Range('D6'); $c.NumberFormat = "@"; $c.Value = '0.6789'; $c.Style = "Normal"
$c = $ws.Range('E6'); $c.NumberFormat = "@"; $c.Value = '  -2.55%  '; $c.Style = "Normal"
$c = $ws.Range('E7'); $c.NumberFormat = "@"; $c.Value = '  -0.03%  '; $c.Style = "Normal"
$c = $ws.Range('D8'); $c.NumberFormat = "@"; $c.Value = '0.2986'; $c.Style = "Normal"
$c = $ws.Range('E8'); $c.NumberFormat = "@"; $c.Value = '  -2.44%  '; $c.Style = "Normal"
$c = $ws.Range('E9'); $c.NumberFormat = "@"; $c.Value = '  -3.20%  '; $c.Style = "Normal"
$c = $ws.Range('E10'); $c.NumberFormat = "@"; $c.Value = '  -2.40%  '; $c.Style = "Normal"
$c = $ws.Range('D11'); $c.NumberFormat = "@"; $c.Value = '0.07665'; $c.Style = "Normal"
$c = $ws.Range('E11'); $c.NumberFormat = "@"; $c.Value = '  -1.08%  '; $c.Style = "Normal"
$c = $ws.Range('D12'); $c.NumberFormat = "@"; $c.Value = '1.833.18'; $c.Style = "Normal"
$c = $ws.Range('E12'); $c.NumberFormat = "@"; $c.Value = '  -1.37%  '; $c.Style = "Normal"
$c = $ws.Range('D13'); $c.NumberFormat = "@"; $c.Value = '5.025'; $c.Style = "Normal"
$c = $ws.Range('E13'); $c.NumberFormat = "@"; $c.Value = '  -2.58%  '; $c.Style = "Normal"
$c = $ws.Range('D14'); $c.NumberFormat = "@"; $c.Value = '0.6762'; $c.Style = "Normal"
$c = $ws.Range('E14'); $c.NumberFormat = "@"; $c.Value = '  -2.36%  '; $c.Style = "Normal"
$c = $ws.Range('D15'); $c.NumberFormat = "@"; $c.Value = '86.07'; $c.Style = "Normal"
$c = $ws.Range('E15'); $c.NumberFormat = "@"; $c.Value = '  -6.65%  '; $c.Style = "Normal"
$c = $ws.Range('D16'); $c.NumberFormat = "@"; $c.Value = '6.159'; $c.Style = "Normal"
$c = $ws.Range('E16'); $c.NumberFormat = "@"; $c.Value = '  -6.24%  '; $c.Style = "Normal"
$c = $ws.Range('D17'); $c.NumberFormat = "@"; $c.Value = '29.110.48'; $c.Style = "Normal"
$c = $ws.Range('E17'); $c.NumberFormat = "@"; $c.Value = '  -1.68%  '; $c.Style = "Normal"
$c = $ws.Range('D18'); $c.NumberFormat = "@"; $c.Value = '0.000008270'; $c.Style = "Normal"
$c = $ws.Range('E18'); $c.NumberFormat = "@"; $c.Value = '  -0.49%  '; $c.Style = "Normal"
$c = $ws.Range('D19'); $c.NumberFormat = "@"; $c.Value = '2.069.87'; $c.Style = "Normal"
$c = $ws.Range('E19'); $c.NumberFormat = "@"; $c.Value = '  -1.18%  '; $c.Style = "Normal"
$c = $ws.Range('D20'); $c.NumberFormat = "@"; $c.Value = '228.11'; $c.Style = "Normal"
$c = $ws.Range('E20'); $c.NumberFormat = "@"; $c.Value = '  -5.36%  '; $c.Style = "Normal"
$c = $ws.Range('D21'); $c.NumberFormat = "@"; $c.Value = '12.49'; $c.Style = "Normal"
$c = $ws.Range('E21'); $c.NumberFormat = "@"; $c.Value = '  -2.18%  '; $c.Style = "Normal"
$c = $ws.Range('D22'); $c.NumberFormat = "@"; $c.Value = '0.9993'; $c.Style = "Normal"
$c = $ws.Range('E22'); $c.NumberFormat = "@"; $c.Value = '  -0.04%  '; $c.Style = "Normal"
$c = $ws.Range('D23'); $c.NumberFormat = "@"; $c.Value = '7.355'; $c.Style = "Normal"
$c = $ws.Range('E23'); $c.NumberFormat = "@"; $c.Value = '  -3.32%  '; $c.Style = "Normal"
$c = $ws.Range('D24'); $c.NumberFormat = "@"; $c.Value = '0.9996'; $c.Style = "Normal"
$c = $ws.Range('E24'); $c.NumberFormat = "@"; $c.Value = '  -0.06%  '; $c.Style = "Normal"
$c = $ws.Range('D25'); $c.NumberFormat = "@"; $c.Value = '160.31'; $c.Style = "Normal"
$c = $ws.Range('E25'); $c.NumberFormat = "@"; $c.Value = '  +0.56%  '; $c.Style = "Normal"
$c = $ws.Range('D26'); $c.NumberFormat = "@"; $c.Value = '0.1436'; $c.Style = "Normal"
$c = $ws.Range('E26'); $c.NumberFormat = "@"; $c.Value = '  -4.28%  '; $c.Style = "Normal"
$c = $ws.Range('D27'); $c.NumberFormat = "@"; $c.Value = '8.707'; $c.Style = "Normal"
$c = $ws.Range('E27'); $c.NumberFormat = "@"; $c.Value = '  -2.43%  '; $c.Style = "Normal"
$c = $ws.Range('E28'); $c.NumberFormat = "@"; $c.Value = '  -1.49%  '; $c.Style = "Normal"
$c = $ws.Range('D29'); $c.NumberFormat = "@"; $c.Value = '1.510'; $c.Style = "Normal"
$c = $ws.Range('E29'); $c.NumberFormat = "@"; $c.Value = '  -1.52%  '; $c.Style = "Normal"
$c = $ws.Range('D30'); $c.NumberFormat = "@"; $c.Value = '4.250'; $c.Style = "Normal"
$c = $ws.Range('E30'); $c.NumberFormat = "@"; $c.Value = '  -0.18%  '; $c.Style = "Normal"
$c = $ws.Range('D31'); $c.NumberFormat = "@"; $c.Value = '4.137'; $c.Style = "Normal"
$c = $ws.Range('E31'); $c.NumberFormat = "@"; $c.Value = '  -1.24%  '; $c.Style = "Normal"
$c = $ws.Range('D32'); $c.NumberFormat = "@"; $c.Value = '1.197'; $c.Style = "Normal"
$c = $ws.Range('E32'); $c.NumberFormat = "@"; $c.Value = '  -0.23%  '; $c.Style = "Normal"
$c = $ws.Range('D33'); $c.NumberFormat = "@"; $c.Value = '0.05414'; $c.Style = "Normal"
$c = $ws.Range('E33'); $c.NumberFormat = "@"; $c.Value = '  +6.32%  '; $c.Style = "Normal"
$c = $ws.Range('D34'); $c.NumberFormat = "@"; $c.Value = '1.864'; $c.Style = "Normal"
$c = $ws.Range('E34'); $c.NumberFormat = "@"; $c.Value = '  -1.82%  '; $c.Style = "Normal"
$c = $ws.Range('D35'); $c.NumberFormat = "@"; $c.Value = '0.7499'; $c.Style = "Normal"
$c = $ws.Range('E35'); $c.NumberFormat = "@"; $c.Value = '  -3.35%  '; $c.Style = "Normal"
$c = $ws.Range('E36'); $c.NumberFormat = "@"; $c.Value = '  -2.12%  '; $c.Style = "Normal"
$c = $ws.Range('D37'); $c.NumberFormat = "@"; $c.Value = '2.681'; $c.Style = "Normal"
$c = $ws.Range('E37'); $c.NumberFormat = "@"; $c.Value = '  -0.10%  '; $c.Style = "Normal"
$c = $ws.Range('D38'); $c.NumberFormat = "@"; $c.Value = '1.305.04'; $c.Style = "Normal"
$c = $ws.Range('E38'); $c.NumberFormat = "@"; $c.Value = '  -1.62%  '; $c.Style = "Normal"
$c = $ws.Range('D39'); $c.NumberFormat = "@"; $c.Value = '0.01815'; $c.Style = "Normal"
$c = $ws.Range('E39'); $c.NumberFormat = "@"; $c.Value = '  -3.14%  '; $c.Style = "Normal"
$c = $ws.Range('D40'); $c.NumberFormat = "@"; $c.Value = '2.715'; $c.Style = "Normal"
$c = $ws.Range('E40'); $c.NumberFormat = "@"; $c.Value = '  -0.62%  '; $c.Style = "Normal"
$c = $ws.Range('D41'); $c.NumberFormat = "@"; $c.Value = '0.9356'; $c.Style = "Normal"
$c = $ws.Range('E41'); $c.NumberFormat = "@"; $c.Value = '  -3.92%  '; $c.Style = "Normal"
$c = $ws.Range('D42'); $c.NumberFormat = "@"; $c.Value = '6.095'; $c.Style = "Normal"
$c = $ws.Range('E42'); $c.NumberFormat = "@"; $c.Value = '  +4.59%  '; $c.Style = "Normal"
$c = $ws.Range('D43'); $c.NumberFormat = "@"; $c.Value = '104.81'; $c.Style = "Normal"
$c = $ws.Range('E43'); $c.NumberFormat = "@"; $c.Value = '  -1.80%  '; $c.Style = "Normal"
$c = $ws.Range('D44'); $c.NumberFormat = "@"; $c.Value = '0.9984'; $c.Style = "Normal"
$c = $ws.Range('E44'); $c.NumberFormat = "@"; $c.Value = '  -0.12%  '; $c.Style = "Normal"

# Rows 45-51: coin list reordered + values updated (Coin, Link, Price, Volume(1h))
$ws.Range('B45').Value = 'BabyDogeCoin'
$ws.Range('C45').Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$c = $ws.Range('D45'); $c.NumberFormat = "@"; $c.Value = '0.00000000126'; $c.Style = "Normal"
$c = $ws.Range('E45'); $c.NumberFormat = "@"; $c.Value = '  -0.09%  '; $c.Style = "Normal"
$ws.Range('B46').Value = 'XinFinNetwork'
$ws.Range('C46').Value = 'https://coinranking.com/coin/77jGXSqWJ1ofG+xinfinnetwork-xdc'
$c = $ws.Range('D46'); $c.NumberFormat = "@"; $c.Value = '0.07874'; $c.Style = "Normal"
$c = $ws.Range('E46'); $c.NumberFormat = "@"; $c.Value = '  +23.19%  '; $c.Style = "Normal"
$ws.Range('B47').Value = 'RocketPoolETH'
$ws.Range('C47').Value = 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth'
$c = $ws.Range('D47'); $c.NumberFormat = "@"; $c.Value = '1.978.66'; $c.Style = "Normal"
$c = $ws.Range('E47'); $c.NumberFormat = "@"; $c.Value = '  -1.14%  '; $c.Style = "Normal"
$ws.Range('B48').Value = 'Mantle'
$ws.Range('C48').Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$c = $ws.Range('D48'); $c.NumberFormat = "@"; $c.Value = '0.5175'; $c.Style = "Normal"
$c = $ws.Range('E48'); $c.NumberFormat = "@"; $c.Value = '  -0.79%  '; $c.Style = "Normal"
$ws.Range('B49').Value = 'Aave'
$ws.Range('C49').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$c = $ws.Range('D49'); $c.NumberFormat = "@"; $c.Value = '63.81'; $c.Style = "Normal"
$c = $ws.Range('E49'); $c.NumberFormat = "@"; $c.Value = '  +0.45%  '; $c.Style = "Normal"
$ws.Range('B50').Value = 'RenderToken'
$ws.Range('C50').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$c = $ws.Range('D50'); $c.NumberFormat = "@"; $c.Value = '1.763'; $c.Style = "Normal"
$c = $ws.Range('E50'); $c.NumberFormat = "@"; $c.Value = '  -0.87%  '; $c.Style = "Normal"
$ws.Range('B51').Value = 'EnergySwap'
$ws.Range('C51').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$c = $ws.Range('D51'); $c.NumberFormat = "@"; $c.Value = '9.396'; $c.Style = "Normal"
$c = $ws.Range('E51'); $c.NumberFormat = "@"; $c.Value = '  -3.97%  '; $c.Style = "Normal"
